# Generate Report for Handoff
# Updates the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" values
# for the 7c18c628-8f4a-42ae-b251-affa754aa880 file row across all report sheets.

$wb = $excel.ActiveWorkbook

# Overview sheet: column G is "Latest HO Xliff Generate Date", row 6 is the
# 7c18c628-8f4a-42ae-b251-affa754aa880.md row.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Cells.Item(6, 7).Value = "2016-08-24 04:39:36"

# zh-cn sheet: column H is "Latest Handoff Datetime", row 6 is the
# 7c18c628-8f4a-42ae-b251-affa754aa880.md row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Cells.Item(6, 8).Value = "2016-08-24 04:39:31"

# de-de sheet: column H is "Latest Handoff Datetime", row 6 is the
# 7c18c628-8f4a-42ae-b251-affa754aa880.md row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Cells.Item(6, 8).Value = "2016-08-24 04:39:36"
